$d = $word.ActiveDocument
$newStyle = $d.Styles.Add("Footnote Block Text", 1)
$newStyle.BaseStyle = $d.Styles.Item("Footnote Text")
$newStyle.NextParagraphStyle = $d.Styles.Item("Footnote Text")
$newStyle.Priority = 9
$newStyle.UnhideWhenUsed = $true
$newStyle.QuickStyle = $true
$newStyle.ParagraphFormat.SpaceBefore = 5
$newStyle.ParagraphFormat.SpaceAfter = 5
$newStyle.ParagraphFormat.FirstLineIndent = 0
$newStyle.ParagraphFormat.LeftIndent = 24
$newStyle.ParagraphFormat.RightIndent = 24
